$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.254.54"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.447.74"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.27"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.11"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "2.443.21"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.46"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000177"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "2.895.64"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "62.192.56"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "2.445.76"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.80"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.54"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.75"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.23"
$ws.Range("E26").Value = "  +3.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "592.11"
$ws.Range("E27").Value = "  -4.40%  "
$ws.Range("D28").Value = "0.0₃0973"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.90"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.136"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.89"
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.377"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.37"
$ws.Range("E39").Value = "  +4.65%  "
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.27"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.88"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.51"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "0.0₆0260"
$ws.Range("E48").Value = "  +18.35%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0521"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.86"
$ws.Range("E51").Value = "  -0.44%  "
